# Remove the "2 x Keys and or Fobs" checklist item row.
# This was a single row in the "Driver Comfort" category (column B)
# with column C holding the item text. Deleting the entire row shifts
# all subsequent rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Cells.Find("2 x Keys and or Fobs")
if ($target -ne $null) {
    $target.EntireRow.Delete()
}
